# no-op, just to test if date fields auto refresh
$p = $ppt.ActivePresentation
Write-Output $p.Slides.Count
